$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This edit adds a new "Foot Length (m)" parameter (and its derived formula)
# to the small lookup table held in columns J:K of the "data" sheet. Doing so
# shifts the existing J:K lookup rows down (the first shift makes room above
# the "Femur COM" entry, the second shift makes room for the new "Foot
# Length" entry just below "Tibia Length"). Columns A:H are not affected by
# the shift - only the J:K helper table moves.
# ---------------------------------------------------------------------------

# Clear out the whole helper area first (J1:K13) since we're going to
# rewrite it from scratch in its final, shifted layout.
$ws.Range("J1:K13").ClearContents()

# --- New layout of the J:K lookup table -----------------------------------
# Row 2: Femur COM (m)
$ws.Range("J2").Value = "Femur COM (m)"
$ws.Range("K2").Value = 0.1950308

# Row 3: Tibia COM (m)
$ws.Range("J3").Value = "Tibia COM (m)"
$ws.Range("K3").Value = 0.18455724409999999

# Row 4: Pelvis Height (m)
$ws.Range("J4").Value = "Pelvis Height (m)"
$ws.Range("K4").Formula = "=1.045"

# Row 6: Femur Length (m)
$ws.Range("J6").Value = "Femur Length (m)"
$ws.Range("K6").Formula = "=K2/(K2+K3)*K4"

# Row 7: Tibia Length (m)
$ws.Range("J7").Value = "Tibia Length (m)"
$ws.Range("K7").Formula = "=K3/(K2+K3)*K4"

# Row 8: Foot Length (m)  -- NEW, computed from the ankle offset (row 23)
$ws.Range("J8").Value = "Foot Length (m)"
$ws.Range("K8").Formula = "=-H23/SIN(RADIANS(F23-G23))"

# Row 10: Height Adjustment Factor
$ws.Range("J10").Value = "Height Adjustment Factor"
$ws.Range("K10").Value = 0.0047697474025435449

# Row 12: Start Time (s)
$ws.Range("J12").Value = "Start Time (s)"
$ws.Range("K12").Value = 0.75

# Row 13: End Time (s)
$ws.Range("J13").Value = "End Time (s)"
$ws.Range("K13").Value = 1.35

# --- Update the H column formulas so they reference the new K rows --------
# (Femur Length now K6, Tibia Length now K7, Height Adjustment Factor now K10)
for ($r = 2; $r -le 74; $r++) {
    $formula = "=B$r-`$K`$6*COS(RADIANS(C$r+D$r))-`$K`$7*COS(RADIANS(C$r+D$r+E$r))-`$K`$10"
    $ws.Range("H$r").Formula = $formula
}

# --- Freeze the header row and update the active selection ----------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("L13").Select()
